$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '61.984.77'
$ws.Range("E2").Value = '  -0.85%  '

$ws.Range("D3").Value = '3.419.87'
$ws.Range("E3").Value = '  -0.67%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").Value = '''410.54'
$ws.Range("E5").Value = '  +0.63%  '

$ws.Range("D6").Value = '''129.87'
$ws.Range("E6").Value = '  -3.04%  '

$ws.Range("E7").Value = '  +8.07%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").Value = '''0.742'
$ws.Range("E9").Value = '  +8.22%  '

$ws.Range("E10").Value = '  +15.10%  '

$ws.Range("D11").Value = '''42.86'
$ws.Range("E11").Value = '  +1.54%  '

$ws.Range("D12").Value = '''0.0000216'
$ws.Range("E12").Value = '  +62.18%  '

$ws.Range("E13").Value = '  +7.48%  '

$ws.Range("E14").Value = '  -0.31%  '

$ws.Range("B15").Value = 'Chainlink'
$ws.Range("C15").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D15").Value = '''21.28'
$ws.Range("E15").Value = '  +6.64%  '

$ws.Range("B16").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C16").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D16").Value = '3.950.97'
$ws.Range("E16").Value = '  -0.81%  '

$ws.Range("D17").Value = '3.408.85'
$ws.Range("E17").Value = '  +1.03%  '

$ws.Range("D18").Value = '''12.19'
$ws.Range("E18").Value = '  +7.03%  '

$ws.Range("D19").Value = '''1.09'
$ws.Range("E19").Value = '  +6.09%  '

$ws.Range("D20").Value = '61.965.74'
$ws.Range("E20").Value = '  -0.86%  '

$ws.Range("D21").Value = '''445.84'
$ws.Range("E21").Value = '  +41.53%  '

$ws.Range("D22").Value = '''91.45'
$ws.Range("E22").Value = '  +8.68%  '

$ws.Range("D23").Value = '''3.17'
$ws.Range("E23").Value = '  -0.41%  '

$ws.Range("D24").Value = '''13.07'
$ws.Range("E24").Value = '  +1.04%  '

$ws.Range("E25").Value = '  +3.50%  '

$ws.Range("D26").Value = '''33.62'
$ws.Range("E26").Value = '  +12.85%  '

$ws.Range("D27").Value = '''8.82'
$ws.Range("E27").Value = '  +6.65%  '

$ws.Range("D28").Value = '''4.74'
$ws.Range("E28").Value = '  +0.19%  '

$ws.Range("D29").Value = '''7.63'
$ws.Range("E29").Value = '  +0.93%  '

$ws.Range("D30").Value = '''2.75'
$ws.Range("E30").Value = '  -1.09%  '

$ws.Range("D31").Value = '''12.03'
$ws.Range("E31").Value = '  +5.48%  '

$ws.Range("E32").Value = '  -0.36%  '

$ws.Range("D33").Value = '''0.169'
$ws.Range("E33").Value = '  -2.45%  '

$ws.Range("D34").Value = '''42.95'
$ws.Range("E34").Value = '  +1.42%  '

$ws.Range("E35").Value = '  -0.05%  '

$ws.Range("D36").Value = '''0.0500'
$ws.Range("E36").Value = '  +2.95%  '

$ws.Range("D37").Value = '''53.64'
$ws.Range("E37").Value = '  +4.19%  '

$ws.Range("D38").Value = '''0.998'
$ws.Range("E38").Value = '  -0.02%  '

$ws.Range("D39").Value = '''3.39'
$ws.Range("E39").Value = '  -1.02%  '

$ws.Range("E40").Value = '  +7.58%  '

$ws.Range("E41").Value = '  -1.32%  '

$ws.Range("D42").Value = '''0.316'
$ws.Range("E42").Value = '  -2.02%  '

$ws.Range("D43").Value = '''141.49'
$ws.Range("E43").Value = '  +2.38%  '

$ws.Range("E44").Value = '  +4.79%  '

$ws.Range("E45").Value = '  -0.42%  '

$ws.Range("E46").Value = '  +7.64%  '

$ws.Range("E47").Value = '  -1.04%  '

$ws.Range("D48").Value = '''22.34'
$ws.Range("E48").Value = '  +4.42%  '

$ws.Range("D49").Value = '3.762.72'
$ws.Range("E49").Value = '  -0.63%  '

$ws.Range("D50").Value = '2.107.24'
$ws.Range("E50").Value = '  -0.99%  '

$ws.Range("E51").Value = '  +13.87%  '
